$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57 (shifts existing rows 57-68 down to 58-69,
# inheriting the formatting - including the date number format in column D -
# from the row above, just like Excel's native row insert).
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly price record.
$ws.Range("A57").Value = 1
$ws.Range("B57").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C57").Value = "Arica y Parinacota"
$ws.Range("D57").Value = 45015
$ws.Range("E57").Value = 15
$ws.Range("F57").Value = 100112028
$ws.Range("G57").Value = "Sandia"
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Segunda"
$ws.Range("J57").Value = 900
$ws.Range("K57").Value = 430
$ws.Range("L57").Value = 450
$ws.Range("M57").Value = 440
$ws.Range("N57").Value = "$/kilo (volumen en unidades)"
$ws.Range("O57").Value = "Perú"
$ws.Range("P57").Value = 440
$ws.Range("Q57").Value = 1
$ws.Range("R57").Value = "Hortaliza"
